$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel
# are temporarily formatted as Text, written, then restored to the default "Normal" style
# so the saved file keeps them as plain text cells (matching the source data) without
# leaving any custom number format applied to the cell.

$ws.Range("D2").Value = "61.634.50"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "2.890.62"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.21%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "2.890.19"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.374.15"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "61.614.59"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "2.891.09"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -11.57%  "
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.60%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  -5.30%  "
$ws.Range("D45").Value = "2.686.82"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "346.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.73%  "
